# Append new daily rows (138-161) of disk_savvy monitoring data to Arkusz1,
# matching the pattern of the existing date/time/files/disk_space columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$newRows = @(
    @(138, 45157, 0.46447916666666672, 77983, 1680),
    @(139, 45158, 0.46090277777777783, 77983, 1680),
    @(140, 45159, 0.37771990740740741, 77983, 1680),
    @(141, 45160, 0.46681712962962968, 78041, 1680),
    @(142, 45161, 0.46094907407407404, 78041, 1680),
    @(143, 45162, 0.46054398148148151, 78178, 1680),
    @(144, 45163, 0.4636805555555556, 78180, 1680),
    @(145, 45164, 0.47665509259259259, 78180, 1680),
    @(146, 45165, 0.46444444444444444, 78187, 1680),
    @(147, 45166, 0.46377314814814818, 78191, 1680),
    @(148, 45167, 0.37917824074074075, 78212, 1680),
    @(149, 45168, 0.46304398148148151, 78242, 1680),
    @(150, 45169, 0.46733796296296298, 78306, 1680),
    @(151, 45170, 0.46157407407407408, 78316, 1680),
    @(152, 45171, 0.46012731481481484, 78325, 1680),
    @(153, 45172, 0.46033564814814815, 78325, 1680),
    @(154, 45173, 0.46032407407407411, 78368, 1690),
    @(155, 45174, 0.46274305555555556, 78368, 1690),
    @(156, 45175, 0.39006944444444441, 78369, 1690),
    @(157, 45176, 0.51321759259259259, 78374, 1690),
    @(158, 45177, 0.46203703703703702, 79071, 1690),
    @(159, 45178, 0.44128472222222226, 79115, 1690),
    @(160, 45179, 0.4636805555555556, 79115, 1690),
    @(161, 45180, 0.46097222222222217, 79119, 1690)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Match date/time number formats used by the existing rows (reuses the
# workbook's existing styles instead of creating new ones).
$ws.Range("A138:A161").NumberFormat = "m/d/yy"
$ws.Range("B138:B161").NumberFormat = "h:mm:ss"

# Column widths for the newly populated B:D columns, best-fit to content
# (mirrors column A, which was already best-fit/custom width). Target
# best-fit widths are B=8.140625, C=6, D=10.5703125 characters; the inputs
# below are chosen so the host's width-rounding lands on those values.
$ws.Columns.Item(2).ColumnWidth = 7.3
$ws.Columns.Item(3).ColumnWidth = 5.15
$ws.Columns.Item(4).ColumnWidth = 9.65

# Move selection to the next empty row below the appended data, same as
# Excel leaves it after typing/pasting a block of rows.
$ws.Range("A162").Select()

Write-Host "Appended $($newRows.Count) rows (138-161) to Arkusz1"
